$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 340.4375
$ws.Range("I2").Value = 334.46155
$ws.Range("J2").Value = 366.33334
$ws.Range("K2").Value = 334.46155
$ws.Range("L2").Value = 366.33334
$ws.Range("M2").Value = -221.46155
$ws.Range("N2").Value = -592.33334

$ws.Range("H3").Value = 41574.8
$ws.Range("J3").Value = 41574.8
$ws.Range("L3").Value = 41574.8
$ws.Range("N3").Value = -41802.8

$ws.Range("H9").Value = 5183.8887
$ws.Range("I9").Value = 9234.700000000001
$ws.Range("J9").Value = 120.375
$ws.Range("K9").Value = 9234.700000000001
$ws.Range("L9").Value = 120.375
$ws.Range("M9").Value = -9065.700000000001
$ws.Range("N9").Value = -458.375

$ws.Range("H29").Value = 3
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 9
$ws.Range("M29").Value = 272

$ws.Range("H38").Value = 550.0833
$ws.Range("I38").Value = 550.0833
$ws.Range("K38").Value = 1650.2499
$ws.Range("M38").Value = -1278.2499

$ws.Range("H51").Value = 2599.9333
$ws.Range("I51").Value = 2499.9285
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 2499.9285
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = -2015.9285
$ws.Range("N51").Value = -4968

$ws.Range("H58").Value = 571.6667
$ws.Range("I58").Value = 571.6667
$ws.Range("K58").Value = 1715.0001
$ws.Range("M58").Value = -1565.0001

$ws.Range("H86").Value = 6184.4287
$ws.Range("J86").Value = 7672
$ws.Range("L86").Value = 7672
$ws.Range("N86").Value = -9918

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H89").Value = 6184.4287
$ws.Range("J89").Value = 7672
$ws.Range("L89").Value = 38360
$ws.Range("N89").Value = -49592

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H102").Value = 41574.8
$ws.Range("J102").Value = 41574.8
$ws.Range("L102").Value = 41574.8
$ws.Range("N102").Value = -48064.8

$ws.Range("H135").Value = 1234.4667
$ws.Range("I135").Value = 1036.9286
$ws.Range("K135").Value = 9332.357399999999
$ws.Range("M135").Value = -6797.357399999999

$ws.Range("H137").Value = 2213.6155
$ws.Range("I137").Value = 2347.625
$ws.Range("K137").Value = 7042.875
$ws.Range("M137").Value = -4492.875

$ws.Range("H138").Value = 3089.311
$ws.Range("J138").Value = 3839.4312
$ws.Range("L138").Value = 11518.2936
$ws.Range("N138").Value = -21798.2936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9342.6
$ws.Range("I45").Value = 10074.692
$ws.Range("K45").Value = 10074.692
$ws.Range("M45").Value = -9697.691999999999

$ws.Range("H61").Value = 2982.7073
$ws.Range("I61").Value = 2390.8823
$ws.Range("K61").Value = 2390.8823
$ws.Range("M61").Value = -2178.8823

$ws.Range("H63").Value = 4209.9165
$ws.Range("I63").Value = 4138.091
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 4138.091
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -3452.091
$ws.Range("N63").Value = -6372

$ws.Range("H66").Value = 4209.9165
$ws.Range("I66").Value = 4138.091
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 20690.455
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -17258.455
$ws.Range("N66").Value = -31864

$ws.Range("H74").Value = 7656.1377
$ws.Range("I74").Value = 1450.8334
$ws.Range("J74").Value = 37441.6
$ws.Range("K74").Value = 1450.8334
$ws.Range("L74").Value = 37441.6
$ws.Range("M74").Value = -576.8334
$ws.Range("N74").Value = -39189.6

$ws.Range("H77").Value = 7656.1377
$ws.Range("I77").Value = 1450.8334
$ws.Range("J77").Value = 37441.6
$ws.Range("K77").Value = 7254.166999999999
$ws.Range("L77").Value = 187208
$ws.Range("M77").Value = -2886.166999999999
$ws.Range("N77").Value = -195944

$ws.Range("H110").Value = 9743.125
$ws.Range("I110").Value = 11535.454
$ws.Range("K110").Value = 11535.454
$ws.Range("M110").Value = -9490.454

$ws.Range("H122").Value = 2151.52
$ws.Range("I122").Value = 2029.8
$ws.Range("K122").Value = 6089.4
$ws.Range("M122").Value = -3639.4

$ws.Range("H136").Value = 2982.7073
$ws.Range("I136").Value = 2390.8823
$ws.Range("K136").Value = 7172.646900000001
$ws.Range("M136").Value = -4622.646900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 8969.125
$ws.Range("I11").Value = 23382.666
$ws.Range("J11").Value = 321
$ws.Range("K11").Value = 23382.666
$ws.Range("L11").Value = 321
$ws.Range("M11").Value = -23242.666
$ws.Range("N11").Value = -601

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 443.9
$ws.Range("I22").Value = 442.5
$ws.Range("J22").Value = 449.5
$ws.Range("K22").Value = 442.5
$ws.Range("L22").Value = 449.5
$ws.Range("M22").Value = -92.5
$ws.Range("N22").Value = -1149.5

$ws.Range("H28").Value = 18007.375
$ws.Range("J28").Value = 17294.285
$ws.Range("L28").Value = 17294.285
$ws.Range("N28").Value = -17784.285

$ws.Range("H31").Value = 34211.79
$ws.Range("I31").Value = 49182
$ws.Range("J31").Value = 8013.9165
$ws.Range("K31").Value = 49182
$ws.Range("L31").Value = 8013.9165
$ws.Range("M31").Value = -48887
$ws.Range("N31").Value = -8603.916499999999

$ws.Range("H34").Value = 34211.79
$ws.Range("I34").Value = 49182
$ws.Range("J34").Value = 8013.9165
$ws.Range("K34").Value = 49182
$ws.Range("L34").Value = 8013.9165
$ws.Range("M34").Value = -48980
$ws.Range("N34").Value = -8417.916499999999

$ws.Range("H58").Value = 2372.162
$ws.Range("I58").Value = 2781.8
$ws.Range("J58").Value = 1518.75
$ws.Range("K58").Value = 2781.8
$ws.Range("L58").Value = 1518.75
$ws.Range("M58").Value = -2578.8
$ws.Range("N58").Value = -1924.75

$ws.Range("H99").Value = 10832.667
$ws.Range("J99").Value = 10832.667
$ws.Range("L99").Value = 10832.667
$ws.Range("N99").Value = -13828.667

$ws.Range("H126").Value = 10832.667
$ws.Range("J126").Value = 10832.667
$ws.Range("L126").Value = 32498.001
$ws.Range("N126").Value = -37438.001

$ws.Range("H134").Value = 16244.077
$ws.Range("I134").Value = 6061.1816
$ws.Range("K134").Value = 18183.5448
$ws.Range("M134").Value = -15648.5448

$ws.Range("H136").Value = 2372.162
$ws.Range("I136").Value = 2781.8
$ws.Range("J136").Value = 1518.75
$ws.Range("K136").Value = 8345.400000000001
$ws.Range("L136").Value = 4556.25
$ws.Range("M136").Value = -5795.400000000001
$ws.Range("N136").Value = -9656.25

$ws.Range("H141").Value = 69983.336
$ws.Range("J141").Value = 69983.336
$ws.Range("L141").Value = 69983.336
$ws.Range("N141").Value = -80343.336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8085.75
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 23514.5
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 23514.5
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = -24054.5

$ws.Range("H73").Value = 8085.75
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 23514.5
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 23514.5
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = -25386.5

$ws.Range("H80").Value = 3098.75
$ws.Range("I80").Value = 2965
$ws.Range("K80").Value = 2965
$ws.Range("M80").Value = -1967

$ws.Range("H83").Value = 3098.75
$ws.Range("I83").Value = 2965
$ws.Range("K83").Value = 14825
$ws.Range("M83").Value = -9833

$ws.Range("H122").Value = 2896.6155
$ws.Range("I122").Value = 1831.8334
$ws.Range("K122").Value = 5495.5002
$ws.Range("M122").Value = -3045.5002

$ws.Range("H126").Value = 12457.917
$ws.Range("I126").Value = 19832.5
$ws.Range("J126").Value = 5083.3335
$ws.Range("K126").Value = 59497.5
$ws.Range("L126").Value = 15250.0005
$ws.Range("M126").Value = -57027.5
$ws.Range("N126").Value = -20190.0005

$ws.Range("H132").Value = 4508.1665
$ws.Range("I132").Value = 4337.375
$ws.Range("J132").Value = 4849.75
$ws.Range("K132").Value = 13012.125
$ws.Range("L132").Value = 14549.25
$ws.Range("M132").Value = -10482.125
$ws.Range("N132").Value = -19609.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5171.2
$ws.Range("I16").Value = 1464
$ws.Range("K16").Value = 1464
$ws.Range("M16").Value = -1294

$ws.Range("H22").Value = 1518.3478
$ws.Range("I22").Value = 1894
$ws.Range("J22").Value = 1482.5714
$ws.Range("K22").Value = 1894
$ws.Range("L22").Value = 1482.5714
$ws.Range("M22").Value = -1599
$ws.Range("N22").Value = -2072.5714

$ws.Range("H27").Value = 1518.3478
$ws.Range("I27").Value = 1894
$ws.Range("J27").Value = 1482.5714
$ws.Range("K27").Value = 1894
$ws.Range("L27").Value = 1482.5714
$ws.Range("M27").Value = -1787
$ws.Range("N27").Value = -1696.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 693.5
$ws.Range("I14").Value = 577.8
$ws.Range("J14").Value = 1106.7142
$ws.Range("K14").Value = 577.8
$ws.Range("L14").Value = 1106.7142
$ws.Range("M14").Value = -409.8
$ws.Range("N14").Value = -1442.7142

$ws.Range("H19").Value = 3666.3333
$ws.Range("J19").Value = 2999
$ws.Range("L19").Value = 2999
$ws.Range("N19").Value = -3347

$ws.Range("H132").Value = 6399.4
$ws.Range("I132").Value = 6399.4
$ws.Range("K132").Value = 19198.2
$ws.Range("M132").Value = -16668.2
